$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price list with latest scraped values.
# Price cells (column D) are forced to Text format before assignment so that
# values such as "1.002" / "0.4530" / "27.304.45" (which use "." as a grouping
# separator or carry significant trailing zeros) are preserved exactly as text
# instead of being auto-converted to numeric values by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.304.45"
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.855.49"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.77"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4530"
$ws.Range("E7").Value = "  -4.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3866"
$ws.Range("E8").Value = "  -4.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.28"
$ws.Range("E9").Value = "  -9.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07926"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.015"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.38"
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.853.72"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.908"
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.144"
$ws.Range("E15").Value = "  -4.94%  "
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.89"
$ws.Range("E17").Value = "  -5.24%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001025"
$ws.Range("E18").Value = "  -3.76%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06553"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.06"
$ws.Range("E20").Value = "  -6.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.502"
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.286.29"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("E24").Value = "  -4.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.287"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.064.58"
$ws.Range("E26").Value = "  -3.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.68"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.90"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.471"
$ws.Range("E30").Value = "  -4.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.34"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09308"
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9353"
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.464"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.583"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.271"
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02224"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06002"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.073"
$ws.Range("E40").Value = "  -11.08%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5911"
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1885"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.14"
$ws.Range("E44").Value = "  -8.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.285"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5624"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.97"
$ws.Range("E47").Value = "  -6.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.371"
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.914"
$ws.Range("E49").Value = "  -5.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06742"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.72"
$ws.Range("E51").Value = "  -1.00%  "
